$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the text number format on the header row and the "Row" (A) column,
# exactly as the upstream CSV -> XLSX export step does on every run.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A20").NumberFormat = "@"

# Refresh the prediction-score column (B) with the latest values produced by
# the quadratic-svm scoring run (outputs-r202, from ful-path.csv).
$ws.Range("B2").Value = 285716.37492780597
$ws.Range("B3").Value = -50347.130308536754
$ws.Range("B4").Value = 296501.35840191541
$ws.Range("B5").Value = 24940.444782851904
$ws.Range("B6").Value = 294066.39692613785
$ws.Range("B7").Value = 308236.81110559159
$ws.Range("B8").Value = -74547.322424563346
$ws.Range("B9").Value = 166118.23430032656
$ws.Range("B10").Value = 282873.40981155506
$ws.Range("B11").Value = 319705.58578585571
$ws.Range("B12").Value = 229580.96301123279
$ws.Range("B13").Value = 329823.46698166232
$ws.Range("B14").Value = 279222.83020849217
$ws.Range("B15").Value = 136400.74829315906
$ws.Range("B16").Value = 303175.01493531937
$ws.Range("B17").Value = 31656.899797970749
$ws.Range("B18").Value = 310390.59049913165
$ws.Range("B19").Value = 94639.060718508088
$ws.Range("B20").Value = 220744.92762393283
